$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to stay text while we write numeric-looking strings,
# then strip the temporary formatting so cells end up back at the default style.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "52.126.39"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "2.904.41"
$ws.Range("E3").Value = "  +3.42%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "354.22"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").Value = "114.03"
$ws.Range("E6").Value = "  +1.53%  "
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.624"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "39.56"
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("D11").Value = "0.0878"
$ws.Range("E11").Value = "  +4.65%  "
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").Value = "19.79"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").Value = "3.364.20"
$ws.Range("E15").Value = "  +3.71%  "
$ws.Range("D16").Value = "2.895.68"
$ws.Range("E16").Value = "  +3.09%  "
$ws.Range("D17").Value = "0.986"
$ws.Range("E17").Value = "  +3.22%  "
$ws.Range("D18").Value = "52.157.02"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").Value = "3.35"
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("D20").Value = "7.61"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").Value = "14.08"
$ws.Range("E21").Value = "  +3.95%  "
$ws.Range("D22").Value = "0.0₃0981"
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").Value = "71.04"
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("D24").Value = "269.85"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("D25").Value = "2.80"
$ws.Range("E25").Value = "  +1.40%  "
$ws.Range("E26").Value = "  +11.75%  "
$ws.Range("D27").Value = "26.80"
$ws.Range("E27").Value = "  +2.44%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "10.66"
$ws.Range("E29").Value = "  +2.30%  "
$ws.Range("E30").Value = "  +13.51%  "
$ws.Range("D31").Value = "6.79"
$ws.Range("E31").Value = "  +10.92%  "
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("D33").Value = "37.36"
$ws.Range("E33").Value = "  -4.93%  "
$ws.Range("D34").Value = "6.10"
$ws.Range("E34").Value = "  +10.50%  "
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "3.33"
$ws.Range("E38").Value = "  +4.82%  "
$ws.Range("D39").Value = "18.80"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("D41").Value = "2.73"
$ws.Range("E41").Value = "  +8.68%  "
$ws.Range("E42").Value = "  +1.00%  "
$ws.Range("E43").Value = "  +5.31%  "
$ws.Range("E44").Value = "  -2.19%  "
$ws.Range("D45").Value = "118.42"
$ws.Range("E45").Value = "  -1.44%  "
$ws.Range("D46").Value = "2.51"
$ws.Range("E46").Value = "  +0.87%  "
$ws.Range("D47").Value = "3.54"
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("D48").Value = "2.181.53"
$ws.Range("E48").Value = "  +2.86%  "
$ws.Range("D49").Value = "0.256"
$ws.Range("E49").Value = "  +16.12%  "
$ws.Range("E50").Value = "  +11.92%  "
$ws.Range("D51").Value = "0.955"
$ws.Range("E51").Value = "  -2.26%  "

$ws.Range("D2:D51").ClearFormats()
